$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 2
    7  = 1
    8  = 1
    9  = 2
    10 = 2
    11 = 0
    12 = 0
    13 = 2
    14 = 3
    15 = 1
    16 = 0
    17 = 2
    18 = 1
    19 = 3
    20 = 0
    21 = 1
    22 = 2
    23 = 0
    24 = 1
    25 = 0
    26 = 1
    27 = 3
    28 = 1
    29 = 3
    30 = 0
    31 = 0
    32 = 1
    34 = 1
    35 = 1
    36 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
